$wb = $excel.ActiveWorkbook

# --- "paper" sheet: specific consumption projections ---
$wsPaper = $wb.Worksheets.Item("paper")

# Pattern A rows: Electricity (B) becomes blank, Heat (C) gets a computed value, Hydrogen (D) becomes blank
$wsPaper.Range("B2").Value = ""
$wsPaper.Range("C2").Value = 0.8700847
$wsPaper.Range("D2").Value = ""
$wsPaper.Range("B3").Value = ""
$wsPaper.Range("C3").Value = -5.0731063
$wsPaper.Range("D3").Value = ""
$wsPaper.Range("B4").Value = ""
$wsPaper.Range("C4").Value = 1.7631324
$wsPaper.Range("D4").Value = ""
$wsPaper.Range("B6").Value = ""
$wsPaper.Range("C6").Value = 1.2620384
$wsPaper.Range("D6").Value = ""
$wsPaper.Range("B11").Value = ""
$wsPaper.Range("C11").Value = -0.2147037
$wsPaper.Range("D11").Value = ""
$wsPaper.Range("B12").Value = ""
$wsPaper.Range("C12").Value = 4.2137493
$wsPaper.Range("D12").Value = ""
$wsPaper.Range("B13").Value = ""
$wsPaper.Range("C13").Value = -0.1167503
$wsPaper.Range("D13").Value = ""
$wsPaper.Range("B14").Value = ""
$wsPaper.Range("C14").Value = 0.09166059999999999
$wsPaper.Range("D14").Value = ""
$wsPaper.Range("B15").Value = ""
$wsPaper.Range("C15").Value = 2.3427368
$wsPaper.Range("D15").Value = ""
$wsPaper.Range("B16").Value = ""
$wsPaper.Range("C16").Value = 1.9525423
$wsPaper.Range("D16").Value = ""
$wsPaper.Range("B18").Value = ""
$wsPaper.Range("C18").Value = -5.249878
$wsPaper.Range("D18").Value = ""
$wsPaper.Range("B19").Value = ""
$wsPaper.Range("C19").Value = 2.4530077
$wsPaper.Range("D19").Value = ""
$wsPaper.Range("B21").Value = ""
$wsPaper.Range("C21").Value = 0.0659781
$wsPaper.Range("D21").Value = ""
$wsPaper.Range("B22").Value = ""
$wsPaper.Range("C22").Value = 2.6669077
$wsPaper.Range("D22").Value = ""
$wsPaper.Range("B23").Value = ""
$wsPaper.Range("C23").Value = 1.8430915
$wsPaper.Range("D23").Value = ""
$wsPaper.Range("B26").Value = ""
$wsPaper.Range("C26").Value = 0.0270267
$wsPaper.Range("D26").Value = ""
$wsPaper.Range("B29").Value = ""
$wsPaper.Range("C29").Value = -0.2755971
$wsPaper.Range("D29").Value = ""
$wsPaper.Range("B31").Value = ""
$wsPaper.Range("C31").Value = 1.8578062
$wsPaper.Range("D31").Value = ""
$wsPaper.Range("B34").Value = ""
$wsPaper.Range("C34").Value = -0.7921222999999999
$wsPaper.Range("D34").Value = ""
$wsPaper.Range("B35").Value = ""
$wsPaper.Range("C35").Value = 0.4206817
$wsPaper.Range("D35").Value = ""

# Pattern B rows: Electricity (B) gets a computed value, Heat (C) and Hydrogen (D) become blank
$wsPaper.Range("B7").Value = 14.8435663
$wsPaper.Range("C7").Value = ""
$wsPaper.Range("D7").Value = ""
$wsPaper.Range("B8").Value = 3.0854156
$wsPaper.Range("C8").Value = ""
$wsPaper.Range("D8").Value = ""
$wsPaper.Range("B9").Value = 2.0036631
$wsPaper.Range("C9").Value = ""
$wsPaper.Range("D9").Value = ""
$wsPaper.Range("B10").Value = 3.8029549
$wsPaper.Range("C10").Value = ""
$wsPaper.Range("D10").Value = ""
$wsPaper.Range("B24").Value = 6.2365862
$wsPaper.Range("C24").Value = ""
$wsPaper.Range("D24").Value = ""
$wsPaper.Range("B33").Value = -13.9284813
$wsPaper.Range("C33").Value = ""
$wsPaper.Range("D33").Value = ""

# Pattern C rows: Electricity (B) and Heat (C) get the EU-aggregate values, Hydrogen (D) stays 0
$wsPaper.Range("B5").Value = 2.9946903
$wsPaper.Range("C5").Value = 7.4707965
$wsPaper.Range("B17").Value = 2.9946903
$wsPaper.Range("C17").Value = 7.4707965
$wsPaper.Range("B20").Value = 2.9946903
$wsPaper.Range("C20").Value = 7.4707965
$wsPaper.Range("B25").Value = 2.9946903
$wsPaper.Range("C25").Value = 7.4707965
$wsPaper.Range("B27").Value = 2.9946903
$wsPaper.Range("C27").Value = 7.4707965
$wsPaper.Range("B28").Value = 2.9946903
$wsPaper.Range("C28").Value = 7.4707965
$wsPaper.Range("B30").Value = 2.9946903
$wsPaper.Range("C30").Value = 7.4707965
$wsPaper.Range("B32").Value = 2.9946903
$wsPaper.Range("C32").Value = 7.4707965

# --- "cement" sheet: every row gets the aggregate Electricity/Heat values, max subst. pct reset to 0 ---
$wsCement = $wb.Worksheets.Item("cement")
for ($r = 2; $r -le 35; $r++) {
    $wsCement.Range("B$r").Value = 0.42092
    $wsCement.Range("C$r").Value = 2.97908
    $wsCement.Range("E$r").Value = 0
}

# --- "glass" sheet: every row gets the aggregate Electricity/Heat values, max subst. pct reset to 0 ---
$wsGlass = $wb.Worksheets.Item("glass")
for ($r = 2; $r -le 35; $r++) {
    $wsGlass.Range("B$r").Value = 1.3898961
    $wsGlass.Range("C$r").Value = 6.1331626
    $wsGlass.Range("E$r").Value = 0
}

# --- alu_prim, alu_sec, copper_prim, copper_sec: reset max subst. of heat with H2 pct from 20 to 0 ---
foreach ($name in @("alu_prim", "alu_sec", "copper_prim", "copper_sec")) {
    $ws = $wb.Worksheets.Item($name)
    for ($r = 2; $r -le 35; $r++) {
        $ws.Range("E$r").Value = 0
    }
}

